$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bl_vendor")

# Remember column B's width before inserting, so the new column C can match it.
$regNoWidth = $ws.Columns("B:B").ColumnWidth

# Insert a new column before C (Company Name, Address lines, Phone/Fax shift right by one).
$ws.Columns("C:C").Insert()

# New header in the freshly inserted column.
$ws.Range("C1").Value = "New Reg No"
$ws.Range("C1").Font.Bold = $true
$ws.Columns("C:C").ColumnWidth = $regNoWidth

# Three additional headers appended after the existing last column (I1 = Fax No).
$ws.Range("J1").Value = "Contact Name"
$ws.Range("J1").Font.Bold = $true

$ws.Range("L1").Value = "Tin No"
$ws.Range("L1").Font.Bold = $true

$ws.Range("K1").Value = "IC No"
$ws.Range("K1").Font.Bold = $true

$ws.Range("M1").Select()
